# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" figures and the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refresh the "last updated" timestamp -------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 13:52"

# --- Finlandia (row 58 / sheet row 54) -------------------------------------
$ws.Range("B54").Value = 4740
$ws.Range("C54").Value = 45
$ws.Range("E54").Value = 2047

# --- Ranks 97-100 (sheet rows 93-96): Republica de Chipre, Albania,
#     Principado de Andorra, Senegal -----------------------------------------
$ws.Range("B93").Value = 823
$ws.Range("C93").Value = 87
$ws.Range("D93").Value = 296
$ws.Range("E93").Value = 518
$ws.Range("F93").Value = 1
$ws.Range("H93").Value = 9

$ws.Range("B94").Value = 822
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 148
$ws.Range("E94").Value = 659
$ws.Range("F94").Value = 15
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 15

$ws.Range("B95").Value = 750
$ws.Range("C95").Value = 14
$ws.Range("D95").Value = 431
$ws.Range("E95").Value = 289
$ws.Range("F95").Value = 4
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 30

$ws.Range("B96").Value = 743
$ws.Range("D96").Value = 385
$ws.Range("E96").Value = 318
$ws.Range("F96").Value = 17
$ws.Range("H96").Value = 40

# --- Ranks 115-120 (sheet rows 111-116): Malta, Jordania, Taiwan,
#     Reunion, Mali, Mayotte -------------------------------------------------
$ws.Range("B111").Value = 460
$ws.Range("C111").Value = 59
$ws.Range("D111").Value = 235
$ws.Range("E111").Value = 221
$ws.Range("F111").Value = 4

$ws.Range("B112").Value = 458
$ws.Range("C112").Value = 8
$ws.Range("D112").Value = 303
$ws.Range("E112").Value = 151
$ws.Range("F112").Value = 1
$ws.Range("H112").Value = 4

$ws.Range("B113").Value = 449
$ws.Range("D113").Value = 342
$ws.Range("E113").Value = 100
$ws.Range("F113").Value = 5
$ws.Range("H113").Value = 7

$ws.Range("B114").Value = 429
$ws.Range("D114").Value = 307
$ws.Range("E114").Value = 116
$ws.Range("F114").Value = 0
$ws.Range("H114").Value = 6

$ws.Range("B115").Value = 418
$ws.Range("D115").Value = 300
$ws.Range("E115").Value = 118
$ws.Range("F115").Value = 2
$ws.Range("H115").Value = 0

$ws.Range("B116").Value = 408
$ws.Range("D116").Value = 113
$ws.Range("E116").Value = 272
$ws.Range("F116").Value = 0
$ws.Range("H116").Value = 23

# --- Madagascar (row 150 / sheet row 141) -----------------------------------
$ws.Range("D141").Value = 82
$ws.Range("E141").Value = 46

# --- Ranks 150-151 (sheet rows 146-147): Bermudas, Cabo Verde ---------------
$ws.Range("B146").Value = 114
$ws.Range("C146").Value = 5
$ws.Range("D146").Value = 2
$ws.Range("E146").Value = 111
$ws.Range("F146").Value = 0
$ws.Range("H146").Value = 1

$ws.Range("B147").Value = 110
$ws.Range("D147").Value = 44
$ws.Range("E147").Value = 60
$ws.Range("F147").Value = 10
$ws.Range("H147").Value = 6
